$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.088.88"
Set-TextValue $ws.Range("E2") "  +4.21%  "
Set-TextValue $ws.Range("D3") "1.909.84"
Set-TextValue $ws.Range("E3") "  +5.39%  "
Set-TextValue $ws.Range("D4") "0.9995"
Set-TextValue $ws.Range("E4") "  +0.04%  "
Set-TextValue $ws.Range("D5") "251.78"
Set-TextValue $ws.Range("E5") "  +1.08%  "
Set-TextValue $ws.Range("D6") "0.9995"
Set-TextValue $ws.Range("E6") "  +0.03%  "
Set-TextValue $ws.Range("D7") "0.5110"
Set-TextValue $ws.Range("E7") "  +2.94%  "
Set-TextValue $ws.Range("B8") "Cardano"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D8") "0.3024"
Set-TextValue $ws.Range("E8") "  +8.62%  "
Set-TextValue $ws.Range("B9") "OKB"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D9") "44.89"
Set-TextValue $ws.Range("E9") "  +3.41%  "
Set-TextValue $ws.Range("D10") "0.06824"
Set-TextValue $ws.Range("E10") "  +5.84%  "
Set-TextValue $ws.Range("D11") "1.911.50"
Set-TextValue $ws.Range("E11") "  +5.52%  "
Set-TextValue $ws.Range("D12") "17.29"
Set-TextValue $ws.Range("E12") "  +2.85%  "
Set-TextValue $ws.Range("E13") "  +3.76%  "
Set-TextValue $ws.Range("D14") "0.7083"
Set-TextValue $ws.Range("E14") "  +8.92%  "
Set-TextValue $ws.Range("D15") "86.85"
Set-TextValue $ws.Range("E15") "  +2.73%  "
Set-TextValue $ws.Range("D16") "4.892"
Set-TextValue $ws.Range("E16") "  +3.78%  "
Set-TextValue $ws.Range("D17") "30.069.87"
Set-TextValue $ws.Range("E17") "  +4.21%  "
Set-TextValue $ws.Range("D18") "0.000008197"
Set-TextValue $ws.Range("E18") "  +10.65%  "
Set-TextValue $ws.Range("D19") "1.000"
Set-TextValue $ws.Range("E19") "  +0.16%  "
Set-TextValue $ws.Range("D20") "13.05"
Set-TextValue $ws.Range("E20") "  +6.25%  "
Set-TextValue $ws.Range("D21") "2.157.00"
Set-TextValue $ws.Range("E21") "  +5.40%  "
Set-TextValue $ws.Range("D22") "0.9987"
Set-TextValue $ws.Range("E22") "  +0.01%  "
Set-TextValue $ws.Range("D23") "4.834"
Set-TextValue $ws.Range("E23") "  +5.19%  "
Set-TextValue $ws.Range("D24") "5.745"
Set-TextValue $ws.Range("E24") "  +7.15%  "
Set-TextValue $ws.Range("D25") "9.246"
Set-TextValue $ws.Range("E25") "  +3.08%  "
Set-TextValue $ws.Range("D26") "147.18"
Set-TextValue $ws.Range("E26") "  +2.80%  "
Set-TextValue $ws.Range("D27") "135.27"
Set-TextValue $ws.Range("E27") "  +2.04%  "
Set-TextValue $ws.Range("D28") "17.11"
Set-TextValue $ws.Range("E28") "  +3.20%  "
Set-TextValue $ws.Range("D29") "2.013"
Set-TextValue $ws.Range("E29") "  +6.28%  "
Set-TextValue $ws.Range("E30") "  -0.56%  "
Set-TextValue $ws.Range("D31") "4.256"
Set-TextValue $ws.Range("E31") "  +1.95%  "
Set-TextValue $ws.Range("D32") "0.08806"
Set-TextValue $ws.Range("E32") "  +5.12%  "
Set-TextValue $ws.Range("D33") "4.028"
Set-TextValue $ws.Range("E33") "  +4.93%  "
Set-TextValue $ws.Range("D34") "0.05081"
Set-TextValue $ws.Range("E34") "  +2.11%  "
Set-TextValue $ws.Range("D35") "1.145"
Set-TextValue $ws.Range("E35") "  +4.63%  "
Set-TextValue $ws.Range("D36") "0.7163"
Set-TextValue $ws.Range("E36") "  +5.81%  "
Set-TextValue $ws.Range("D37") "2.689"
Set-TextValue $ws.Range("E37") "  -0.71%  "
Set-TextValue $ws.Range("D38") "2.812"
Set-TextValue $ws.Range("E38") "  +1.63%  "
Set-TextValue $ws.Range("E39") "  -1.03%  "
Set-TextValue $ws.Range("D40") "0.9660"
Set-TextValue $ws.Range("E40") "  +0.34%  "
Set-TextValue $ws.Range("D41") "0.01696"
Set-TextValue $ws.Range("E41") "  +6.07%  "
Set-TextValue $ws.Range("D42") "6.183"
Set-TextValue $ws.Range("E42") "  +2.24%  "
Set-TextValue $ws.Range("D43") "0.4333"
Set-TextValue $ws.Range("E43") "  +5.40%  "
Set-TextValue $ws.Range("D44") "105.58"
Set-TextValue $ws.Range("E44") "  +5.61%  "
Set-TextValue $ws.Range("D45") "0.9990"
Set-TextValue $ws.Range("E45") "  -0.01%  "
Set-TextValue $ws.Range("D46") "7.681"
Set-TextValue $ws.Range("E46") "  +6.03%  "
Set-TextValue $ws.Range("D47") "0.1280"
Set-TextValue $ws.Range("E47") "  +4.58%  "
Set-TextValue $ws.Range("D48") "0.05740"
Set-TextValue $ws.Range("E48") "  +3.89%  "
Set-TextValue $ws.Range("D49") "33.21"
Set-TextValue $ws.Range("E49") "  +4.99%  "
Set-TextValue $ws.Range("D50") "8.431"
Set-TextValue $ws.Range("E50") "  +3.62%  "
Set-TextValue $ws.Range("D51") "0.3821"
Set-TextValue $ws.Range("E51") "  +4.87%  "
